$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 25.00614441910898
$ws.Range("C2").Value = 12.30731053561313
$ws.Range("D2").Value = 4.060659819108497
$ws.Range("E2").Value = 9.760591166727075
$ws.Range("F2").Value = 53.8300619044932
$ws.Range("I2").Value = 37.5689805815659
$ws.Range("J2").Value = 9.653495557233981
$ws.Range("L2").Value = 12.67856184045885

# Row 3
$ws.Range("B3").Value = 24.74140692267896
$ws.Range("C3").Value = 11.95175431647375
$ws.Range("D3").Value = 4.034675934778845
$ws.Range("E3").Value = 9.769613835736275
$ws.Range("F3").Value = 53.70344327964792
$ws.Range("I3").Value = 37.55785225907291
$ws.Range("J3").Value = 9.670375217775872
$ws.Range("L3").Value = 12.68597716097457

# Row 4
$ws.Range("B4").Value = 24.58469654994497
$ws.Range("C4").Value = 11.73228476847314
$ws.Range("D4").Value = 4.018360675747212
$ws.Range("E4").Value = 9.775533184945582
$ws.Range("F4").Value = 53.63809740843654
$ws.Range("I4").Value = 37.5585900983377
$ws.Range("J4").Value = 9.68131401363596
$ws.Range("L4").Value = 12.69273076665195

# Row 5
$ws.Range("B5").Value = 24.52237692755122
$ws.Range("C5").Value = 11.64271556975237
$ws.Range("D5").Value = 4.01162133703517
$ws.Range("E5").Value = 9.778040985024994
$ws.Range("F5").Value = 53.61459065031933
$ws.Range("I5").Value = 37.56078876182038
$ws.Range("J5").Value = 9.685916554616913
$ws.Range("L5").Value = 12.69603650654158

# Row 6
$ws.Range("B6").Value = 24.51212394373675
$ws.Range("C6").Value = 11.62783935030672
$ws.Range("D6").Value = 4.010496795684248
$ws.Range("E6").Value = 9.778463184749697
$ws.Range("F6").Value = 53.61087601526222
$ws.Range("I6").Value = 37.56126829473169
$ws.Range("J6").Value = 9.686689566000677
$ws.Range("L6").Value = 12.69661886155771

# Row 7
$ws.Range("B7").Value = 24.58384975218972
$ws.Range("C7").Value = 11.73107713427289
$ws.Range("D7").Value = 4.018270153979175
$ws.Range("E7").Value = 9.775566618563666
$ws.Range("F7").Value = 53.6377677430311
$ws.Range("I7").Value = 37.5586120739417
$ws.Range("J7").Value = 9.681375497920021
$ws.Range("L7").Value = 12.69277310736494

# Row 8
$ws.Range("B8").Value = 24.91369552549715
$ws.Range("C8").Value = 12.1850466711856
$ws.Range("D8").Value = 4.051775543144227
$ws.Range("E8").Value = 9.763623568833104
$ws.Range("F8").Value = 53.78383412554831
$ws.Range("I8").Value = 37.56356909335588
$ws.Range("J8").Value = 9.659196644569686
$ws.Range("L8").Value = 12.68066210169591

# Row 9
$ws.Range("B9").Value = 25.60351802926553
$ws.Range("C9").Value = 13.05962207465618
$ws.Range("D9").Value = 4.114613666487101
$ws.Range("E9").Value = 9.743204044116906
$ws.Range("F9").Value = 54.16831647701391
$ws.Range("I9").Value = 37.63357610836263
$ws.Range("J9").Value = 9.620244974427532
$ws.Range("L9").Value = 12.67435775374529

# Row 10
$ws.Range("B10").Value = 26.1319158658457
$ws.Range("C10").Value = 13.68462226798604
$ws.Range("D10").Value = 4.159034210943943
$ws.Range("E10").Value = 9.730018175672992
$ws.Range("F10").Value = 54.50993715559498
$ws.Range("I10").Value = 37.72200093700684
$ws.Range("J10").Value = 9.594369800315073
$ws.Range("L10").Value = 12.6803336106378

# Row 11
$ws.Range("B11").Value = 26.37597358210039
$ws.Range("C11").Value = 13.963652488355
$ws.Range("D11").Value = 4.178862930283477
$ws.Range("E11").Value = 9.724411289378757
$ws.Range("F11").Value = 54.67799241578483
$ws.Range("I11").Value = 37.77028858014427
$ws.Range("J11").Value = 9.583188576958653
$ws.Range("L11").Value = 12.68534633564191

# Row 12
$ws.Range("B12").Value = 26.46883680303266
$ws.Range("C12").Value = 14.06844431505551
$ws.Range("D12").Value = 4.186317023389641
$ws.Range("E12").Value = 9.722344190777225
$ws.Range("F12").Value = 54.74342764874428
$ws.Range("I12").Value = 37.78973349258648
$ws.Range("J12").Value = 9.579038897800887
$ws.Range("L12").Value = 12.68757331588326

# Row 13
$ws.Range("B13").Value = 26.44881867511249
$ws.Range("C13").Value = 14.04591590541688
$ws.Range("D13").Value = 4.184714083251881
$ws.Range("E13").Value = 9.722786884873518
$ws.Range("F13").Value = 54.72925550072607
$ws.Range("I13").Value = 37.78549412880533
$ws.Range("J13").Value = 9.579928857196968
$ws.Range("L13").Value = 12.68707909384386

# Row 14
$ws.Range("B14").Value = 26.3836050922645
$ws.Range("C14").Value = 13.9722918259324
$ws.Range("D14").Value = 4.179477278695827
$ws.Range("E14").Value = 9.72424010446902
$ws.Range("F14").Value = 54.68333995098747
$ws.Range("I14").Value = 37.77186509211911
$ws.Range("J14").Value = 9.582845490788698
$ws.Range("L14").Value = 12.68552296880117

# Row 15
$ws.Range("B15").Value = 26.34371513950159
$ws.Range("C15").Value = 13.9270784293374
$ws.Range("D15").Value = 4.176262466475491
$ws.Range("E15").Value = 9.725137545022919
$ws.Range("F15").Value = 54.65544856359605
$ws.Range("I15").Value = 37.76366788634239
$ws.Range("J15").Value = 9.584642994626581
$ws.Range("L15").Value = 12.68461257580518

# Row 16
$ws.Range("B16").Value = 26.11603292349072
$ws.Range("C16").Value = 13.66627091675537
$ws.Range("D16").Value = 4.157730726296216
$ws.Range("E16").Value = 9.730392460022181
$ws.Range("F16").Value = 54.49920702922238
$ws.Range("I16").Value = 37.71900749689097
$ws.Range("J16").Value = 9.59511235098781
$ws.Range("L16").Value = 12.68005209031207

# Row 17
$ws.Range("B17").Value = 25.97724005637663
$ws.Range("C17").Value = 13.50483802759884
$ws.Range("D17").Value = 4.146265047276716
$ws.Range("E17").Value = 9.733716305990074
$ws.Range("F17").Value = 54.40658272081157
$ws.Range("I17").Value = 37.69367523640279
$ws.Range("J17").Value = 9.601685692304565
$ws.Range("L17").Value = 12.67784124289586

# Row 18
$ws.Range("B18").Value = 25.89776225969459
$ws.Range("C18").Value = 13.41149455387271
$ws.Range("D18").Value = 4.139634686927137
$ws.Range("E18").Value = 9.735664948171026
$ws.Range("F18").Value = 54.35449992927575
$ws.Range("I18").Value = 37.6798637667904
$ws.Range("J18").Value = 9.605522014766109
$ws.Range("L18").Value = 12.67678565019492

# Row 19
$ws.Range("B19").Value = 25.87091556736409
$ws.Range("C19").Value = 13.37980918282601
$ws.Range("D19").Value = 4.13738365149646
$ws.Range("E19").Value = 9.736331060165524
$ws.Range("F19").Value = 54.33707100681732
$ws.Range("I19").Value = 37.67531777583554
$ws.Range("J19").Value = 9.606830472762583
$ws.Range("L19").Value = 12.67646537958958

# Row 20
$ws.Range("B20").Value = 25.99197897019615
$ws.Range("C20").Value = 13.52207454336574
$ws.Range("D20").Value = 4.147489274762105
$ws.Range("E20").Value = 9.733358663989586
$ws.Range("F20").Value = 54.41631949744005
$ws.Range("I20").Value = 37.69629334380543
$ws.Range("J20").Value = 9.600980206440536
$ws.Range("L20").Value = 12.67805424099423

# Row 21
$ws.Range("B21").Value = 26.402748541621
$ws.Range("C21").Value = 13.99394146899006
$ws.Range("D21").Value = 4.181016938938222
$ws.Range("E21").Value = 9.723811737319103
$ws.Range("F21").Value = 54.69677790434425
$ws.Range("I21").Value = 37.77583680956084
$ws.Range("J21").Value = 9.581986516905916
$ws.Range("L21").Value = 12.68597112773575

# Row 22
$ws.Range("B22").Value = 26.67375777348134
$ws.Range("C22").Value = 14.29721364333667
$ws.Range("D22").Value = 4.202610782909639
$ws.Range("E22").Value = 9.717899225924919
$ws.Range("F22").Value = 54.89053159704223
$ws.Range("I22").Value = 37.83458029917207
$ws.Range("J22").Value = 9.57006486645065
$ws.Range("L22").Value = 12.69306093623591

# Row 23
$ws.Range("B23").Value = 26.52891030994232
$ws.Range("C23").Value = 14.13585399160777
$ws.Range("D23").Value = 4.191114945756697
$ws.Range("E23").Value = 9.721024986340742
$ws.Range("F23").Value = 54.78617297665271
$ws.Range("I23").Value = 37.80260979872451
$ws.Range("J23").Value = 9.576382793757871
$ws.Range("L23").Value = 12.68910211313433

# Row 24
$ws.Range("B24").Value = 25.98531451621679
$ws.Range("C24").Value = 13.51428356797193
$ws.Range("D24").Value = 4.146935921695703
$ws.Range("E24").Value = 9.733520236438515
$ws.Range("F24").Value = 54.4119138604892
$ws.Range("I24").Value = 37.69510735395497
$ws.Range("J24").Value = 9.601298978341635
$ws.Range("L24").Value = 12.67795727338164

# Row 25
$ws.Range("B25").Value = 25.41277102371339
$ws.Range("C25").Value = 12.82556824623507
$ws.Range("D25").Value = 4.097920975809591
$ws.Range("E25").Value = 9.748408170706938
$ws.Range("F25").Value = 54.05385976446701
$ws.Range("I25").Value = 37.60815429206531
$ws.Range("J25").Value = 9.63029895415117
$ws.Range("L25").Value = 12.6741977592793
